$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, preserving the cell's existing
# style (these price/volume cells are plain inline strings with no
# explicit number format, so we briefly force Text format to stop Excel
# from reinterpreting numeric-looking / percent-looking strings as
# numbers, then restore the original style).
function Set-TextValue($rangeAddr, $text) {
    $cell = $ws.Range($rangeAddr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue "D2" '330.81'
Set-TextValue "E2" '0.34%'
Set-TextValue "D3" '41.52'
Set-TextValue "E3" '0.85%'
Set-TextValue "D4" '5.687'
Set-TextValue "E4" '-0.22%'
Set-TextValue "D5" '0.08430'
Set-TextValue "D6" '8.811'
Set-TextValue "E6" '1.18%'
Set-TextValue "D7" '1.992'
Set-TextValue "E7" '-1.68%'
Set-TextValue "D8" '4.488'
Set-TextValue "E8" '-0.75%'
Set-TextValue "E10" '0.43%'
Set-TextValue "E11" '0.18%'
Set-TextValue "D12" '0.1964'
Set-TextValue "E12" '1.02%'
Set-TextValue "D13" '0.09351'
Set-TextValue "E13" '0.74%'
Set-TextValue "D14" '0.03956'
Set-TextValue "E14" '6.98%'
Set-TextValue "D15" '0.1063'
Set-TextValue "E15" '0.72%'
Set-TextValue "D16" '0.001302'
Set-TextValue "E16" '0.07%'
Set-TextValue "D17" '0.006115'
Set-TextValue "E17" '-3.07%'
Set-TextValue "D18" '3.434'
Set-TextValue "E18" '1.53%'
Set-TextValue "E19" '1.11%'
Set-TextValue "D20" '9.182'
Set-TextValue "E20" '10.85%'
Set-TextValue "D21" '0.1364'
Set-TextValue "E21" '-3.80%'
Set-TextValue "D22" '0.2513'
Set-TextValue "E22" '-5.27%'
Set-TextValue "D23" '0.04417'
Set-TextValue "E23" '-0.43%'
Set-TextValue "D24" '0.001247'
Set-TextValue "E24" '-1.14%'
Set-TextValue "D25" '0.004400'
Set-TextValue "E25" '0.82%'
Set-TextValue "E26" '-3.99%'
Set-TextValue "D27" '0.0003997'
Set-TextValue "E27" '0.09%'
Set-TextValue "D39" '0.02845'
Set-TextValue "E39" '0.71%'
Set-TextValue "D40" '0.05530'
Set-TextValue "E40" '1.20%'
Set-TextValue "D41" '0.007912'
Set-TextValue "E41" '3.78%'
Set-TextValue "D42" '0.1436'
Set-TextValue "E42" '1.41%'
Set-TextValue "D43" '0.008976'
Set-TextValue "E43" '-9.84%'
Set-TextValue "E44" '-1.79%'
Set-TextValue "D45" '0.01099'
Set-TextValue "E45" '-7.67%'
Set-TextValue "E46" '8.09%'
Set-TextValue "E47" '0.00%'
Set-TextValue "D48" '0.003246'
Set-TextValue "E48" '8.22%'
Set-TextValue "D49" '0.002282'
Set-TextValue "E49" '-0.05%'
Set-TextValue "D50" '0.00002104'
Set-TextValue "E50" '0.00%'
Set-TextValue "D51" '0.0002004'
Set-TextValue "E51" '0.00%'
